# Automatische test-sync: 2025-08-14 22:08:50
# Adds the new mail-log entry (row 39) to the "Logs" sheet, extends the
# conditional-formatting ranges to cover it, and bumps the "Intern verzoek /
# Actie voor medewerker" tally on the "Dashboard" sheet from 30 to 31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 39

$ws.Cells.Item($newRow, 1).Value = "Datasheet opvragen"
$ws.Cells.Item($newRow, 2).Value = "retour@testbedrijf123.nl"
$ws.Cells.Item($newRow, 3).Value = "Kun je mij de datasheet van de VentiQ-250 sturen?"
$ws.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar documentatie@testbedrijf123.nl."
$ws.Cells.Item($newRow, 6).Value = "2025-08-14 22:08:40"
$ws.Cells.Item($newRow, 7).Value = "Nee"
$ws.Cells.Item($newRow, 8).Value = "Ja"
$ws.Cells.Item($newRow, 9).Value = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Extend the existing conditional formatting rules so they cover the new row
# (was row range 2:38, now 2:39) for columns D, G, H, I and J.
$cfColumns = "D", "G", "H", "I", "J"
foreach ($col in $cfColumns) {
    $oldRange = $ws.Range("$col`2:$col`38")
    $newRange = $ws.Range("$col`2:$col`39")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for the category that received the new row.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 31
